# Refresh the cryptos price/volume snapshot (GitHub Actions data pull).
# Price (column D) and Volume(1h) (column E) are plain text cells in this
# sheet (e.g. "29.926.56", "  +2.55%  "), not numbers - some of the new
# Price values happen to look like valid numbers (e.g. "246.25"), so for
# those we briefly force Text number-format before writing the value (so
# Excel doesn't auto-convert it to a numeric cell) and then clear the
# formatting again right away so the cell's style stays exactly as it was.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.926.56"
$ws.Range("E2").Value = "  +2.55%  "
$ws.Range("D3").Value = "1.864.97"
$ws.Range("E3").Value = "  +2.19%  "
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "246.25"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +1.99%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6405"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +3.52%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "47.29"
$ws.Range("D8").ClearFormats()
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3004"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +3.67%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07485"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +1.97%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "24.38"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +6.11%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07679"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +0.12%  "
$ws.Range("D13").Value = "1.879.98"
$ws.Range("E13").Value = "  +3.03%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.080"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +2.14%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6906"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +4.34%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "84.33"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +2.62%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000009466"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +6.09%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.107"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +4.51%  "
$ws.Range("D19").Value = "29.917.49"
$ws.Range("E19").Value = "  +2.60%  "
$ws.Range("D20").Value = "2.127.59"
$ws.Range("E20").Value = "  +2.91%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "241.77"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +1.90%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "12.70"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +2.19%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.001"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +0.08%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.490"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +4.19%  "
$ws.Range("E25").Value = "  +0.11%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "159.81"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +0.94%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1426"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +0.90%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.598"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +1.86%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "18.06"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +2.30%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.06161"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +10.57%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.508"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +1.44%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.275"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +5.81%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.171"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +1.81%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.138"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +1.14%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.870"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +2.71%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.165"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +3.02%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.7368"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +0.42%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.607"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -0.45%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.871"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +0.94%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01807"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +2.65%  "
$ws.Range("D41").Value = "1.224.21"
$ws.Range("E41").Value = "  +0.62%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9284"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +1.04%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.295"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -0.27%  "
$ws.Range("D44").Value = "2.039.28"
$ws.Range("E44").Value = "  +3.52%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.003"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +0.26%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "102.36"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +0.83%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "66.82"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +3.21%  "
$ws.Range("E48").Value = "  +5.93%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.5095"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +0.27%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.392"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +3.42%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4102"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +2.34%  "
